$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: format column A as Text (adds the new numFmt=49 "@" style) ---
$ws.Columns.Item(1).NumberFormat = "@"

# --- Step 2: re-enter existing rows 2-13 col A as text (same digits, now text) ---
$ws.Range("A2").Value = '5413458064879'
$ws.Range("A3").Value = '3564700713235'
$ws.Range("A4").Value = '4056489030942'
$ws.Range("A5").Value = '3270160202706'
$ws.Range("A6").Value = '213400032605'
$ws.Range("A7").Value = '5413458017578'
$ws.Range("A8").Value = '3261055947093'
$ws.Range("A9").Value = '3760152390691'
$ws.Range("A10").Value = '3770009392051'
$ws.Range("A11").Value = '3770009392044'
$ws.Range("A12").Value = '2227171006138'
$ws.Range("A13").Value = '213311065143'

# --- Step 3: append new product rows 38-49 (category D) ---
# Row 38
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = '5400141319460'
$ws.Range("B38").Value = 'Kipfilet filet de poulet'
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = 'd'
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = 'Fish Meat Eggs'
$ws.Range("G38").Value = 'Meat'
$ws.Range("H38").Value = 544
$ws.Range("I38").Value = 1.7
$ws.Range("J38").Value = 1.2
$ws.Range("K38").Value = 0.84
$ws.Range("L38").Value = 17
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 5.6
$ws.Range("O38").Value = 1.7
$ws.Range("P38").Value = 1.2
$ws.Range("Q38").Value = 2.1
$ws.Range("R38").Value = 9
$ws.Range("S38").Value = 'n'

# Row 39
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = '3256221978063'
$ws.Range("B39").Value = 'Emincés de poulet traité en salaison cuite et grillé'
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = 'd'
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 'Fish Meat Eggs'
$ws.Range("G39").Value = 'Processed meat'
$ws.Range("H39").Value = 728
$ws.Range("I39").Value = 2.5
$ws.Range("J39").Value = 1
$ws.Range("K39").Value = 0.742
$ws.Range("L39").Value = 21
$ws.Range("M39").Value = 1.5
$ws.Range("N39").Value = 9
$ws.Range("O39").Value = 2.5
$ws.Range("P39").Value = 1
$ws.Range("Q39").Value = 1.855
$ws.Range("R39").Value = 33
$ws.Range("S39").Value = 'n'

# Row 40
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = '5400112591703'
$ws.Range("B40").Value = 'Lamelles de poulet'
$ws.Range("C40").Value = 3
$ws.Range("D40").Value = 'd'
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 'Fish Meat Eggs'
$ws.Range("G40").Value = 'Meat'
$ws.Range("H40").Value = 451
$ws.Range("I40").Value = 1.2000000476837
$ws.Range("J40").Value = 1
$ws.Range("K40").Value = 0.83999996185304
$ws.Range("L40").Value = 19
$ws.Range("M40").Value = 0.1
$ws.Range("N40").Value = 3
$ws.Range("O40").Value = 1.2000000476837
$ws.Range("P40").Value = 1
$ws.Range("Q40").Value = 2.0999999046326
$ws.Range("R40").Value = 9
$ws.Range("S40").Value = 'n'

# Row 41
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = '5400141323986'
$ws.Range("B41").Value = 'Filet de poulet'
$ws.Range("C41").Value = 5
$ws.Range("D41").Value = 'd'
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 'Fish Meat Eggs'
$ws.Range("G41").Value = 'Meat'
$ws.Range("H41").Value = 544
$ws.Range("I41").Value = 1.7
$ws.Range("J41").Value = 1.2
$ws.Range("K41").Value = 0.84
$ws.Range("L41").Value = 17
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 5.6
$ws.Range("O41").Value = 1.7
$ws.Range("P41").Value = 1.2
$ws.Range("Q41").Value = 2.1
$ws.Range("R41").Value = 10
$ws.Range("S41").Value = 'n'

# Row 42
$ws.Range("A42").Value = 27063613
$ws.Range("A42").NumberFormat = "@"
$ws.Range("B42").Value = 'Lamelles poulet précuites Aldi Delifin'
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 'd'
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 'Fish Meat Eggs'
$ws.Range("G42").Value = 'Meat'
$ws.Range("H42").Value = 607
$ws.Range("I42").Value = 0.80000001192093
$ws.Range("J42").Value = 2
$ws.Range("K42").Value = 0.91999998092652
$ws.Range("L42").Value = 25
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 4
$ws.Range("O42").Value = 0.80000001192093
$ws.Range("P42").Value = 2
$ws.Range("Q42").Value = 2.2999999523163
$ws.Range("R42").Value = 12
$ws.Range("S42").Value = 'n'

# Row 43
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = '3095759136015'
$ws.Range("B43").Value = 'Emincés de poulet "J''Aime"'
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 'd'
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = 'Fish Meat Eggs'
$ws.Range("G43").Value = 'Meat'
$ws.Range("H43").Value = 698
$ws.Range("I43").Value = 2
$ws.Range("J43").Value = 0.8
$ws.Range("K43").Value = 0.76
$ws.Range("L43").Value = 22
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 8.4
$ws.Range("O43").Value = 2
$ws.Range("P43").Value = 0.8
$ws.Range("Q43").Value = 1.9
$ws.Range("R43").Value = 39
$ws.Range("S43").Value = 'n'

# Row 44
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = '8437002390228'
$ws.Range("B44").Value = 'Roti de Pollo'
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = 'd'
$ws.Range("E44").Value = 4
$ws.Range("F44").Value = 'Fish Meat Eggs'
$ws.Range("G44").Value = 'Meat'
$ws.Range("H44").Value = 602
$ws.Range("I44").Value = 3.1
$ws.Range("J44").Value = 0.7
$ws.Range("K44").Value = 0.84
$ws.Range("L44").Value = 12
$ws.Range("M44").Value = 1
$ws.Range("N44").Value = 8.2
$ws.Range("O44").Value = 3.1
$ws.Range("P44").Value = 0.7
$ws.Range("Q44").Value = 2.1
$ws.Range("R44").Value = 7
$ws.Range("S44").Value = 'n'

# Row 45
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = '7613035287167'
$ws.Range("B45").Value = 'Knacki 100% poulet'
$ws.Range("C45").Value = 2
$ws.Range("D45").Value = 'd'
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 'Fish Meat Eggs'
$ws.Range("G45").Value = 'Processed meat'
$ws.Range("H45").Value = 827
$ws.Range("I45").Value = 4.5
$ws.Range("J45").Value = 1.5
$ws.Range("K45").Value = 0.72
$ws.Range("L45").Value = 14
$ws.Range("M45").Value = 1
$ws.Range("N45").Value = 15
$ws.Range("O45").Value = 4.5
$ws.Range("P45").Value = 1.5
$ws.Range("Q45").Value = 1.8
$ws.Range("R45").Value = 8
$ws.Range("S45").Value = 'n'

# Row 46
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = '8437005669765'
$ws.Range("B46").Value = 'Pollo relleno asado al horno'
$ws.Range("C46").Value = 6
$ws.Range("D46").Value = 'd'
$ws.Range("E46").Value = 4
$ws.Range("F46").Value = 'Fish Meat Eggs'
$ws.Range("G46").Value = 'Meat'
$ws.Range("H46").Value = 695
$ws.Range("I46").Value = 3.3
$ws.Range("J46").Value = 1.3
$ws.Range("K46").Value = 0.6
$ws.Range("L46").Value = 18.5
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 9.5
$ws.Range("O46").Value = 3.3
$ws.Range("P46").Value = 1.3
$ws.Range("Q46").Value = 1.5
$ws.Range("R46").Value = 9
$ws.Range("S46").Value = 'n'

# Row 47
$ws.Range("A47").Value = 27007600
$ws.Range("A47").NumberFormat = "@"
$ws.Range("B47").Value = 'Poulet hawaï'
$ws.Range("C47").Value = 5
$ws.Range("D47").Value = 'd'
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 'Salty snacks'
$ws.Range("G47").Value = 'Salty and fatty products'
$ws.Range("H47").Value = 1117
$ws.Range("I47").Value = 2
$ws.Range("J47").Value = 6.5
$ws.Range("K47").Value = 0.6
$ws.Range("L47").Value = 8.6000003814697
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 22
$ws.Range("O47").Value = 2
$ws.Range("P47").Value = 6.5
$ws.Range("Q47").Value = 1.5
$ws.Range("R47").Value = 27
$ws.Range("S47").Value = 'n'

# Row 48
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = '3512690003393'
$ws.Range("B48").Value = 'Allumettes de poulet fumées'
$ws.Range("C48").Value = 7
$ws.Range("D48").Value = 'd'
$ws.Range("E48").Value = 4
$ws.Range("F48").Value = 'Fish Meat Eggs'
$ws.Range("G48").Value = 'Processed meat'
$ws.Range("H48").Value = 343
$ws.Range("I48").Value = 0.1
$ws.Range("J48").Value = 3.1
$ws.Range("K48").Value = 1.4
$ws.Range("L48").Value = 16.8
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0.3
$ws.Range("O48").Value = 0.1
$ws.Range("P48").Value = 3.1
$ws.Range("Q48").Value = 3.5
$ws.Range("R48").Value = 9
$ws.Range("S48").Value = 'n'

# Row 49
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = '3560071011383'
$ws.Range("B49").Value = 'Blanc de poulet'
$ws.Range("C49").Value = 4
$ws.Range("D49").Value = 'd'
$ws.Range("E49").Value = 4
$ws.Range("F49").Value = 'Fish Meat Eggs'
$ws.Range("G49").Value = 'Processed meat'
$ws.Range("H49").Value = 413
$ws.Range("I49").Value = 0.3
$ws.Range("J49").Value = 1.1
$ws.Range("K49").Value = 0.96
$ws.Range("L49").Value = 21
$ws.Range("M49").Value = 0
$ws.Range("N49").Value = 1
$ws.Range("O49").Value = 0.3
$ws.Range("P49").Value = 1.1
$ws.Range("Q49").Value = 2.4
$ws.Range("R49").Value = 9
$ws.Range("S49").Value = 'n'

# --- Step 4: column width + page setup + selection (best-effort cosmetic match) ---
$ws.Columns.Item(1).ColumnWidth = 34.5
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Range("F51").Select()
